$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D (shifts old D:K data to G:N)
$ws.Range("D1:F1").EntireColumn.Insert()

# Copy number formats from the (now-shifted) old D column (now column G) into the new D:F columns
# so the new cells inherit the same date/number formatting as the rest of the table.
$ws.Range("G5:G102").Copy()
$ws.Range("D5:F102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns D, E, F with the latest three quarters of data
# Row 7
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("F7").Value2 = 43281
# Row 8
$ws.Range("D8").Value2 = 7200
$ws.Range("E8").Value2 = 6800
$ws.Range("F8").Value2 = 6400
# Row 9
$ws.Range("D9").Value2 = "NA"
$ws.Range("E9").Value2 = "NA"
$ws.Range("F9").Value2 = "NA"
# Row 10
$ws.Range("D10").Value2 = "NA"
$ws.Range("E10").Value2 = "NA"
$ws.Range("F10").Value2 = "NA"
# Row 12
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("F12").Value2 = "NA"
# Row 13
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = 0
# Row 14
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("F14").Value2 = 0
# Row 15
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("F15").Value2 = 0
# Row 17
$ws.Range("D17").Value2 = 0
$ws.Range("E17").Value2 = 0
$ws.Range("F17").Value2 = 0
# Row 18
$ws.Range("D18").Value2 = 7200
$ws.Range("E18").Value2 = 6800
$ws.Range("F18").Value2 = 6400
# Row 20
$ws.Range("D20").Value2 = -1900
$ws.Range("E20").Value2 = -4000
$ws.Range("F20").Value2 = -3800
# Row 21
$ws.Range("D21").Value2 = "NA"
$ws.Range("E21").Value2 = "NA"
$ws.Range("F21").Value2 = "NA"
# Row 22
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 0
$ws.Range("F22").Value2 = 0
# Row 23
$ws.Range("D23").Value2 = 5300
$ws.Range("E23").Value2 = 2800
$ws.Range("F23").Value2 = 2600
# Row 24
$ws.Range("D24").Value2 = 900
$ws.Range("E24").Value2 = 600
$ws.Range("F24").Value2 = 700
# Row 25
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 0
# Row 26
$ws.Range("D26").Value2 = 4400
$ws.Range("E26").Value2 = 2200
$ws.Range("F26").Value2 = 1900
# Row 27
$ws.Range("D27").Value2 = 4400
$ws.Range("E27").Value2 = 2200
$ws.Range("F27").Value2 = 1900
# Row 28
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 0
# Row 29
$ws.Range("D29").Value2 = 100
$ws.Range("E29").Value2 = "NA"
$ws.Range("F29").Value2 = "NA"
# Row 30
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("F30").Value2 = 0
# Row 31
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("F31").Value2 = 0
# Row 32
$ws.Range("D32").Value2 = 1900
$ws.Range("E32").Value2 = 4000
$ws.Range("F32").Value2 = 3800
# Row 33
$ws.Range("D33").Value2 = 4500
$ws.Range("E33").Value2 = 2200
$ws.Range("F33").Value2 = 1900
# Row 34
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("F34").Value2 = 0
# Row 35
$ws.Range("D35").Value2 = 4500
$ws.Range("E35").Value2 = 2200
$ws.Range("F35").Value2 = 1900
# Row 38
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("F38").Value2 = 43281
# Row 41
$ws.Range("D41").Value2 = 145800
$ws.Range("E41").Value2 = "NA"
$ws.Range("F41").Value2 = "NA"
# Row 42
$ws.Range("D42").Value2 = 147400
$ws.Range("E42").Value2 = "NA"
$ws.Range("F42").Value2 = "NA"
# Row 43
$ws.Range("D43").Value2 = 0
$ws.Range("E43").Value2 = 0
$ws.Range("F43").Value2 = 0
# Row 44
$ws.Range("D44").Value2 = 0
$ws.Range("E44").Value2 = 0
$ws.Range("F44").Value2 = 0
# Row 45
$ws.Range("D45").Value2 = 0
$ws.Range("E45").Value2 = 0
$ws.Range("F45").Value2 = 0
# Row 46
$ws.Range("D46").Value2 = 0
$ws.Range("E46").Value2 = 0
$ws.Range("F46").Value2 = 0
# Row 47
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("F47").Value2 = 0
# Row 48
$ws.Range("D48").Value2 = 14700
$ws.Range("E48").Value2 = "NA"
$ws.Range("F48").Value2 = "NA"
# Row 49
$ws.Range("D49").Value2 = 0
$ws.Range("E49").Value2 = 0
$ws.Range("F49").Value2 = 0
# Row 50
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("F50").Value2 = 0
# Row 51
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("F51").Value2 = 0
# Row 52
$ws.Range("D52").Value2 = 3500
$ws.Range("E52").Value2 = "NA"
$ws.Range("F52").Value2 = "NA"
# Row 53
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("F53").Value2 = 0
# Row 54
$ws.Range("D54").Value2 = 746900
$ws.Range("E54").Value2 = 742200
$ws.Range("F54").Value2 = 725000
# Row 57
$ws.Range("D57").Value2 = 0
$ws.Range("E57").Value2 = 0
$ws.Range("F57").Value2 = 0
# Row 58
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("F58").Value2 = 0
# Row 59
$ws.Range("D59").Value2 = 12400
$ws.Range("E59").Value2 = "NA"
$ws.Range("F59").Value2 = "NA"
# Row 60
$ws.Range("D60").Value2 = 0
$ws.Range("E60").Value2 = 0
$ws.Range("F60").Value2 = 0
# Row 61
$ws.Range("D61").Value2 = 4000
$ws.Range("E61").Value2 = 0
$ws.Range("F61").Value2 = 0
# Row 62
$ws.Range("D62").Value2 = 0
$ws.Range("E62").Value2 = 0
$ws.Range("F62").Value2 = 0
# Row 63
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("F63").Value2 = 0
# Row 64
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("F64").Value2 = 0
# Row 65
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("F65").Value2 = 0
# Row 66
$ws.Range("D66").Value2 = 669800
$ws.Range("E66").Value2 = 651600
$ws.Range("F66").Value2 = 638600
# Row 68
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("F68").Value2 = 0
# Row 69
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("F69").Value2 = 0
# Row 70
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("F70").Value2 = 0
# Row 71
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("F71").Value2 = 0
# Row 72
$ws.Range("D72").Value2 = 6300
$ws.Range("E72").Value2 = "NA"
$ws.Range("F72").Value2 = "NA"
# Row 73
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("F73").Value2 = 0
# Row 74
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("F74").Value2 = 0
# Row 75
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("F75").Value2 = 0
# Row 76
$ws.Range("D76").Value2 = 77100
$ws.Range("E76").Value2 = "NA"
$ws.Range("F76").Value2 = "NA"
# Row 77
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("F77").Value2 = 0
# Row 80
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("F80").Value2 = 43281
# Row 81
$ws.Range("D81").Value2 = 4500
$ws.Range("E81").Value2 = 2200
$ws.Range("F81").Value2 = 1900
# Row 83
$ws.Range("D83").Value2 = 0
$ws.Range("E83").Value2 = 0
$ws.Range("F83").Value2 = 0
# Row 84
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("F84").Value2 = 0
# Row 85
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("F85").Value2 = 0
# Row 86
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("F86").Value2 = 0
# Row 87
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("F87").Value2 = 0
# Row 88
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("F88").Value2 = 0
# Row 89
$ws.Range("D89").Value2 = 0
$ws.Range("E89").Value2 = 0
$ws.Range("F89").Value2 = 0
# Row 91
$ws.Range("D91").Value2 = 0
$ws.Range("E91").Value2 = 0
$ws.Range("F91").Value2 = 0
# Row 92
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("F92").Value2 = 0
# Row 93
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("F93").Value2 = 0
# Row 94
$ws.Range("D94").Value2 = 0
$ws.Range("E94").Value2 = 0
$ws.Range("F94").Value2 = 0
# Row 96
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = 0
$ws.Range("F96").Value2 = 0
# Row 97
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("F97").Value2 = 0
# Row 98
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("F98").Value2 = 0
# Row 99
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("F99").Value2 = 0
# Row 100
$ws.Range("D100").Value2 = 0
$ws.Range("E100").Value2 = 0
$ws.Range("F100").Value2 = 0
# Row 101
$ws.Range("D101").Value2 = 0
$ws.Range("E101").Value2 = 0
$ws.Range("F101").Value2 = 0
# Row 102
$ws.Range("D102").Value2 = 0
$ws.Range("E102").Value2 = 0
$ws.Range("F102").Value2 = 0
